$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.115, -0.109, -0.321, 66, -19, 22),
    @(0.203, 0.501, -0.024, -139, 4, -169),
    @(0.038, 0.061, 0.866, -100, 55, -135),
    @(-0.239, -0.028, 0.706, 49, 1, -129),
    @(0.126, -0.125, 0.294, 96, 18, 44),
    @(0.12, 0.048, 0.305, 36, 28, 127),
    @(-0.075, 0.01, 0.82, -83, -51, -157),
    @(-0.074, -0.207, 0.189, 68, -50, 41),
    @(0.299, -0.232, -0.12, 88, -66, -4),
    @(-0.223, -0.082, 0.187, 139, -24, -51)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}
